$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
  "대표","기자","한국","대해","기술","지역","사업","대한","국회","금지",
  "무단","북한","원내대표","조사","지난해","위해","서울","배포","민주당","의원",
  "문화","통해","총리","때문","서비스","처리","계획","관련","확대","중단",
  "미국","이후","시장","관리","사람","정부","합의","라며","발언","문제",
  "정준영","사용","모델","지난","뉴스","이번","개발","가장","고시원","미세먼지",
  "내용","위원회","진행","위원장","국민","이상","설명","아베","개정안","개선",
  "시간","운항","연설","르노","올해","국가","면서","대통령","증가"
)

$counts = @(
  71,63,62,57,53,45,44,44,40,39,
  38,38,38,37,37,36,35,35,35,35,
  35,34,34,32,32,31,31,30,29,29,
  29,29,29,28,28,28,28,27,26,25,
  25,24,24,24,24,24,24,23,23,22,
  22,22,22,22,22,22,22,22,21,21,
  21,21,21,21,20,20,20,20,20
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $words[$i]
    $ws.Cells.Item($row, 2).Value = $counts[$i]
}

$wb.Save()
